$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Value = "'21"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'21.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("G3").Value = "'21"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.442"
$ws.Range("D4").Style = "Normal"
$ws.Range("G4").Value = "'21"
$ws.Range("G4").Style = "Normal"
$ws.Range("G5").Value = "'21"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'3.420"
$ws.Range("D6").Style = "Normal"
$ws.Range("G6").Value = "'21"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'6.325"
$ws.Range("D7").Style = "Normal"
$ws.Range("G7").Value = "'21"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.8171"
$ws.Range("D8").Style = "Normal"
$ws.Range("G8").Value = "'21"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'1.036"
$ws.Range("D9").Style = "Normal"
$ws.Range("G9").Value = "'21"
$ws.Range("G9").Style = "Normal"
$ws.Range("G10").Value = "'21"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.07287"
$ws.Range("D11").Style = "Normal"
$ws.Range("G11").Value = "'21"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.03142"
$ws.Range("D12").Style = "Normal"
$ws.Range("G12").Value = "'21"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.03138"
$ws.Range("D13").Style = "Normal"
$ws.Range("G13").Value = "'21"
$ws.Range("G13").Style = "Normal"
$ws.Range("G14").Value = "'21"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.09378"
$ws.Range("D15").Style = "Normal"
$ws.Range("G15").Value = "'21"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.001596"
$ws.Range("D16").Style = "Normal"
$ws.Range("G16").Value = "'21"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'0.04815"
$ws.Range("D17").Style = "Normal"
$ws.Range("G17").Value = "'21"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005840"
$ws.Range("D18").Style = "Normal"
$ws.Range("G18").Value = "'21"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.006222"
$ws.Range("D19").Style = "Normal"
$ws.Range("G19").Value = "'21"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'0.004136"
$ws.Range("D20").Style = "Normal"
$ws.Range("G20").Value = "'21"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.0009924"
$ws.Range("D21").Style = "Normal"
$ws.Range("G21").Value = "'21"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001498"
$ws.Range("D22").Style = "Normal"
$ws.Range("G22").Value = "'21"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'3.744"
$ws.Range("D23").Style = "Normal"
$ws.Range("G23").Value = "'21"
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = "'2.176"
$ws.Range("D24").Style = "Normal"
$ws.Range("G24").Value = "'21"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.3233"
$ws.Range("D25").Style = "Normal"
$ws.Range("G25").Value = "'21"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.1308"
$ws.Range("D26").Style = "Normal"
$ws.Range("G26").Value = "'21"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003992"
$ws.Range("D27").Style = "Normal"
$ws.Range("G27").Value = "'21"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'21"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'21"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'21"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'21"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'21"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'21"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'21"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'21"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'21"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'21"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'21"
$ws.Range("G38").Style = "Normal"
$ws.Range("G39").Value = "'21"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.03859"
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").Value = "'21"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.006696"
$ws.Range("D41").Style = "Normal"
$ws.Range("G41").Value = "'21"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.1069"
$ws.Range("D42").Style = "Normal"
$ws.Range("G42").Value = "'21"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.002726"
$ws.Range("D43").Style = "Normal"
$ws.Range("G43").Value = "'21"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.006571"
$ws.Range("D44").Style = "Normal"
$ws.Range("G44").Value = "'21"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005591"
$ws.Range("D45").Style = "Normal"
$ws.Range("G45").Value = "'21"
$ws.Range("G45").Style = "Normal"
$ws.Range("G46").Value = "'21"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.3893"
$ws.Range("D47").Style = "Normal"
$ws.Range("G47").Value = "'21"
$ws.Range("G47").Style = "Normal"
$ws.Range("G48").Value = "'21"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002097"
$ws.Range("D49").Style = "Normal"
$ws.Range("G49").Value = "'21"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.01009"
$ws.Range("D50").Style = "Normal"
$ws.Range("G50").Value = "'21"
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = "'21"
$ws.Range("G51").Style = "Normal"
